$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.02801400538725
$ws.Range("D2").Value = 1.030658468443009
$ws.Range("E2").Value = 1.036685735900681
$ws.Range("F2").Value = 1.044479839310314
$ws.Range("I2").Value = 1.030351808346303
$ws.Range("J2").Value = 1.033168678404457
$ws.Range("K2").Value = 1.033468841885566
$ws.Range("L2").Value = 1.039478769466597
$ws.Range("M2").Value = 1.047250772005269
$ws.Range("N2").Value = 1.015054268736965

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.028935535647547
$ws.Range("D3").Value = 1.031493186519562
$ws.Range("E3").Value = 1.037539489137516
$ws.Range("F3").Value = 1.045497940334492
$ws.Range("I3").Value = 1.030438544036345
$ws.Range("J3").Value = 1.033730781382564
$ws.Range("K3").Value = 1.034111890175055
$ws.Range("L3").Value = 1.040142058333258
$ws.Range("M3").Value = 1.04807957978765
$ws.Range("N3").Value = 1.015242508426825

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029532221470078
$ws.Range("D4").Value = 1.032033970756056
$ws.Range("E4").Value = 1.038092715070924
$ws.Range("F4").Value = 1.046157848587315
$ws.Range("I4").Value = 1.0304932034776
$ws.Range("J4").Value = 1.034094281568128
$ws.Range("K4").Value = 1.034528006684874
$ws.Range("L4").Value = 1.040571399001371
$ws.Range("M4").Value = 1.048616393656533
$ws.Range("N4").Value = 1.015364180781803

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.029783161479916
$ws.Range("D5").Value = 1.032261474445029
$ws.Range("E5").Value = 1.038325479145094
$ws.Range("F5").Value = 1.046435542960212
$ws.Range("I5").Value = 1.030515831188816
$ws.Range("J5").Value = 1.034247043836143
$ws.Range("K5").Value = 1.034702946020196
$ws.Range("L5").Value = 1.040751928255756
$ws.Range("M5").Value = 1.048842193463608
$ws.Range("N5").Value = 1.015415300010652

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.029825300851488
$ws.Range("D6").Value = 1.032299682539623
$ws.Range("E6").Value = 1.038364572236033
$ws.Range("F6").Value = 1.04648218480069
$ws.Range("I6").Value = 1.0305196098754
$ws.Range("J6").Value = 1.034272690143414
$ws.Range("K6").Value = 1.034732319303228
$ws.Range("L6").Value = 1.04078224189958
$ws.Range("M6").Value = 1.048880113438134
$ws.Range("N6").Value = 1.015423881276206

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.02953557417753
$ws.Range("D7").Value = 1.032037010052429
$ws.Range("E7").Value = 1.038095824539289
$ws.Range("F7").Value = 1.046161558097973
$ws.Range("I7").Value = 1.030493507210525
$ws.Range("J7").Value = 1.034096322995381
$ws.Range("K7").Value = 1.034530344218723
$ws.Range("L7").Value = 1.04057381110816
$ws.Range("M7").Value = 1.048619410321778
$ws.Range("N7").Value = 1.015364863965002

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028325358606722
$ws.Range("D8").Value = 1.0309404266174
$ws.Range("E8").Value = 1.036974101310486
$ws.Range("F8").Value = 1.044823676809475
$ws.Range("I8").Value = 1.03038142395178
$ws.Range("J8").Value = 1.033358688272084
$ws.Range("K8").Value = 1.033686157978513
$ws.Range("L8").Value = 1.039702899655155
$ws.Range("M8").Value = 1.047530763129946
$ws.Range("N8").Value = 1.015117912174211

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026195869618124
$ws.Range("D9").Value = 1.029013267878408
$ws.Range("E9").Value = 1.035003595791914
$ws.Range("F9").Value = 1.04247486165209
$ws.Range("I9").Value = 1.030172727996974
$ws.Range("J9").Value = 1.032057256526075
$ws.Range("K9").Value = 1.032198806999615
$ws.Range("L9").Value = 1.038169427397912
$ws.Range("M9").Value = 1.045616461266404
$ws.Range("N9").Value = 1.014681763200137

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024778333452507
$ws.Range("D10").Value = 1.027732048448473
$ws.Range("E10").Value = 1.033694116216783
$ws.Range("F10").Value = 1.040914912250438
$ws.Range("I10").Value = 1.030026105661032
$ws.Range("J10").Value = 1.031188602383285
$ws.Range("K10").Value = 1.031207449883624
$ws.Range("L10").Value = 1.037147973723157
$ws.Range("M10").Value = 1.044343041844542
$ws.Range("N10").Value = 1.014390355847377

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024165040929421
$ws.Range("D11").Value = 1.027178126281488
$ws.Range("E11").Value = 1.033128107740483
$ws.Range("F11").Value = 1.040240858591146
$ws.Range("I11").Value = 1.029960846391237
$ws.Range("J11").Value = 1.030812232912889
$ws.Range("K11").Value = 1.030778245075269
$ws.Range("L11").Value = 1.036705890632276
$ws.Range("M11").Value = 1.0437923114498
$ws.Range("N11").Value = 1.014264025936787

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023937314372993
$ws.Range("D12").Value = 1.02697250464191
$ws.Range("E12").Value = 1.032918019272255
$ws.Range("F12").Value = 1.039990698833925
$ws.Range("I12").Value = 1.02993634066173
$ws.Range("J12").Value = 1.030672397969774
$ws.Range("K12").Value = 1.030618829310482
$ws.Range("L12").Value = 1.036541714352195
$ws.Range("M12").Value = 1.043587847245329
$ws.Range("N12").Value = 1.014217079434962

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023986158965448
$ws.Range("D13").Value = 1.027016605289691
$ws.Range("E13").Value = 1.032963077049624
$ws.Range("F13").Value = 1.040044349256956
$ws.Range("I13").Value = 1.029941609231445
$ws.Range("J13").Value = 1.030702394588319
$ws.Range("K13").Value = 1.030653024062598
$ws.Range("L13").Value = 1.036576929208529
$ws.Range("M13").Value = 1.04363170089609
$ws.Range("N13").Value = 1.014227150606105

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024146215392698
$ws.Range("D14").Value = 1.027161126878526
$ws.Range("E14").Value = 1.033110738648039
$ws.Range("F14").Value = 1.040220175939276
$ws.Range("I14").Value = 1.029958826153757
$ws.Range("J14").Value = 1.03080067481237
$ws.Range("K14").Value = 1.030765067507872
$ws.Range("L14").Value = 1.036692319089189
$ws.Range("M14").Value = 1.043775408296002
$ws.Range("N14").Value = 1.014260145769853

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024244841751131
$ws.Range("D15").Value = 1.027250188628843
$ws.Range("E15").Value = 1.033201738050484
$ws.Range("F15").Value = 1.040328536926807
$ws.Range("I15").Value = 1.029969398897832
$ws.Range("J15").Value = 1.030861223950932
$ws.Range("K15").Value = 1.030834102531699
$ws.Range("L15").Value = 1.036763419015664
$ws.Range("M15").Value = 1.043863964663018
$ws.Range("N15").Value = 1.014280472282166

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024819046485327
$ws.Range("D16").Value = 1.02776882855427
$ws.Range("E16").Value = 1.033731701618197
$ws.Range("F16").Value = 1.040959676884071
$ws.Range("I16").Value = 1.030030399438537
$ws.Range("J16").Value = 1.031213575885327
$ws.Range("K16").Value = 1.031235936137196
$ws.Range("L16").Value = 1.037177317901042
$ws.Range("M16").Value = 1.044379606190739
$ws.Range("N16").Value = 1.014398736848832

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02517936716961
$ws.Range("D17").Value = 1.028094387363279
$ws.Range("E17").Value = 1.034064403712425
$ws.Range("F17").Value = 1.041355953939725
$ws.Range("I17").Value = 1.030068189728834
$ws.Range("J17").Value = 1.031434534310752
$ws.Range("K17").Value = 1.031488012406413
$ws.Range("L17").Value = 1.037437003527597
$ws.Range("M17").Value = 1.044703234305324
$ws.Range("N17").Value = 1.014472881568358

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.025389585269863
$ws.Range("D18").Value = 1.028284362644497
$ws.Range("E18").Value = 1.03425856007732
$ws.Range("F18").Value = 1.041587232059169
$ws.Range("I18").Value = 1.030090061271683
$ws.Range("J18").Value = 1.031563392738848
$ws.Range("K18").Value = 1.031635049884011
$ws.Range("L18").Value = 1.037588494249797
$ws.Range("M18").Value = 1.044892065545059
$ws.Range("N18").Value = 1.014516114565387

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.02546127251519
$ws.Range("D19").Value = 1.02834915322999
$ws.Range("E19").Value = 1.034324778769708
$ws.Range("F19").Value = 1.041666115022284
$ws.Range("I19").Value = 1.030097489900271
$ws.Range("J19").Value = 1.031607326222055
$ws.Range("K19").Value = 1.031685186789965
$ws.Range("L19").Value = 1.037640152114713
$ws.Range("M19").Value = 1.044956463012973
$ws.Range("N19").Value = 1.014530853452345

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025140703080251
$ws.Range("D20").Value = 1.028059449457841
$ws.Range("E20").Value = 1.034028697897921
$ws.Range("F20").Value = 1.041313423042727
$ws.Range("I20").Value = 1.030064152860976
$ws.Range("J20").Value = 1.031410829922741
$ws.Range("K20").Value = 1.031460966415563
$ws.Range("L20").Value = 1.037409139590558
$ws.Range("M20").Value = 1.044668505398304
$ws.Range("N20").Value = 1.014464928024594

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024099080614378
$ws.Range("D21").Value = 1.027118565280923
$ws.Range("E21").Value = 1.033067251773371
$ws.Range("F21").Value = 1.040168393452256
$ws.Range("I21").Value = 1.029953763525066
$ws.Range("J21").Value = 1.030771734682325
$ws.Range("K21").Value = 1.030732073220709
$ws.Range("L21").Value = 1.036658338724203
$ws.Range("M21").Value = 1.043733087235175
$ws.Range("N21").Value = 1.014250430119356

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023444620645057
$ws.Range("D22").Value = 1.026527745232746
$ws.Range("E22").Value = 1.032463634141565
$ws.Range("F22").Value = 1.039449705406333
$ws.Range("I22").Value = 1.029882821343266
$ws.Range("J22").Value = 1.030369710641661
$ws.Range("K22").Value = 1.030273847603334
$ws.Range("L22").Value = 1.03618647223639
$ws.Range("M22").Value = 1.043145541013813
$ws.Range("N22").Value = 1.014115439974335

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023791519332039
$ws.Range("D23").Value = 1.026840878459371
$ws.Range("E23").Value = 1.032783539253019
$ws.Range("F23").Value = 1.039830577822406
$ws.Range("I23").Value = 1.02992057456109
$ws.Range("J23").Value = 1.030582849779278
$ws.Range("K23").Value = 1.030516755766412
$ws.Range("L23").Value = 1.036436599057056
$ws.Range("M23").Value = 1.043456954239096
$ws.Range("N23").Value = 1.014187012697247

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025158173566175
$ws.Range("D24").Value = 1.028075236139154
$ws.Range("E24").Value = 1.034044831519396
$ws.Range("F24").Value = 1.041332640502954
$ws.Range("I24").Value = 1.03006597747588
$ws.Range("J24").Value = 1.031421540986128
$ws.Range("K24").Value = 1.031473187317274
$ws.Range("L24").Value = 1.037421730041098
$ws.Range("M24").Value = 1.044684197697495
$ws.Range("N24").Value = 1.014468521933081

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.026746023904374
$ws.Range("D25").Value = 1.029510864419448
$ws.Range("E25").Value = 1.035512285906611
$ws.Range("F25").Value = 1.043081047523057
$ws.Range("I25").Value = 1.030228003743206
$ws.Range("J25").Value = 1.032393894266493
$ws.Range("K25").Value = 1.032583290733513
$ws.Range("L25").Value = 1.03856571966844
$ws.Range("M25").Value = 1.046110868991548
$ws.Range("N25").Value = 1.014794632904129

